$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 373, pushing the existing
# rows 373-393 down to 377-397. This makes room for a new week's
# Chirimoya price entries (Especial/Primera/Segunda/Tercera).
$ws.Rows("373:376").Insert()

# Common (fixed) column values shared by every data row in this sheet.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"
$unidad = "$/bandeja 10 kilos"
$origen = "Provincia de Limarí"
$kgPorUnidad = 10

# New row 373: Especial
$r = 373
$ws.Cells.Item($r,1).Value2 = $mercadoId
$ws.Cells.Item($r,2).Value2 = $mercado
$ws.Cells.Item($r,3).Value2 = $region
$ws.Cells.Item($r,4).Value2 = 45265
$ws.Cells.Item($r,5).Value2 = $codreg
$ws.Cells.Item($r,6).Value2 = $tipo
$ws.Cells.Item($r,7).Value2 = $productoId
$ws.Cells.Item($r,8).Value2 = $producto
$ws.Cells.Item($r,9).Value2 = $categoriaId
$ws.Cells.Item($r,10).Value2 = $categoria
$ws.Cells.Item($r,11).Value2 = $variedad
$ws.Cells.Item($r,12).Value2 = "Especial"
$ws.Cells.Item($r,13).Value2 = 440
$ws.Cells.Item($r,14).Value2 = 17000
$ws.Cells.Item($r,15).Value2 = 18000
$ws.Cells.Item($r,16).Value2 = 17500
$ws.Cells.Item($r,17).Value2 = $unidad
$ws.Cells.Item($r,18).Value2 = $origen
$ws.Cells.Item($r,19).Value2 = 1750
$ws.Cells.Item($r,20).Value2 = $kgPorUnidad

# New row 374: Primera
$r = 374
$ws.Cells.Item($r,1).Value2 = $mercadoId
$ws.Cells.Item($r,2).Value2 = $mercado
$ws.Cells.Item($r,3).Value2 = $region
$ws.Cells.Item($r,4).Value2 = 45265
$ws.Cells.Item($r,5).Value2 = $codreg
$ws.Cells.Item($r,6).Value2 = $tipo
$ws.Cells.Item($r,7).Value2 = $productoId
$ws.Cells.Item($r,8).Value2 = $producto
$ws.Cells.Item($r,9).Value2 = $categoriaId
$ws.Cells.Item($r,10).Value2 = $categoria
$ws.Cells.Item($r,11).Value2 = $variedad
$ws.Cells.Item($r,12).Value2 = "Primera"
$ws.Cells.Item($r,13).Value2 = 360
$ws.Cells.Item($r,14).Value2 = 14000
$ws.Cells.Item($r,15).Value2 = 15000
$ws.Cells.Item($r,16).Value2 = 14500
$ws.Cells.Item($r,17).Value2 = $unidad
$ws.Cells.Item($r,18).Value2 = $origen
$ws.Cells.Item($r,19).Value2 = 1450
$ws.Cells.Item($r,20).Value2 = $kgPorUnidad

# New row 375: Segunda
$r = 375
$ws.Cells.Item($r,1).Value2 = $mercadoId
$ws.Cells.Item($r,2).Value2 = $mercado
$ws.Cells.Item($r,3).Value2 = $region
$ws.Cells.Item($r,4).Value2 = 45265
$ws.Cells.Item($r,5).Value2 = $codreg
$ws.Cells.Item($r,6).Value2 = $tipo
$ws.Cells.Item($r,7).Value2 = $productoId
$ws.Cells.Item($r,8).Value2 = $producto
$ws.Cells.Item($r,9).Value2 = $categoriaId
$ws.Cells.Item($r,10).Value2 = $categoria
$ws.Cells.Item($r,11).Value2 = $variedad
$ws.Cells.Item($r,12).Value2 = "Segunda"
$ws.Cells.Item($r,13).Value2 = 300
$ws.Cells.Item($r,14).Value2 = 10000
$ws.Cells.Item($r,15).Value2 = 11000
$ws.Cells.Item($r,16).Value2 = 10500
$ws.Cells.Item($r,17).Value2 = $unidad
$ws.Cells.Item($r,18).Value2 = $origen
$ws.Cells.Item($r,19).Value2 = 1050
$ws.Cells.Item($r,20).Value2 = $kgPorUnidad

# New row 376: Tercera
$r = 376
$ws.Cells.Item($r,1).Value2 = $mercadoId
$ws.Cells.Item($r,2).Value2 = $mercado
$ws.Cells.Item($r,3).Value2 = $region
$ws.Cells.Item($r,4).Value2 = 45265
$ws.Cells.Item($r,5).Value2 = $codreg
$ws.Cells.Item($r,6).Value2 = $tipo
$ws.Cells.Item($r,7).Value2 = $productoId
$ws.Cells.Item($r,8).Value2 = $producto
$ws.Cells.Item($r,9).Value2 = $categoriaId
$ws.Cells.Item($r,10).Value2 = $categoria
$ws.Cells.Item($r,11).Value2 = $variedad
$ws.Cells.Item($r,12).Value2 = "Tercera"
$ws.Cells.Item($r,13).Value2 = 240
$ws.Cells.Item($r,14).Value2 = 7000
$ws.Cells.Item($r,15).Value2 = 8000
$ws.Cells.Item($r,16).Value2 = 7500
$ws.Cells.Item($r,17).Value2 = $unidad
$ws.Cells.Item($r,18).Value2 = $origen
$ws.Cells.Item($r,19).Value2 = 750
$ws.Cells.Item($r,20).Value2 = $kgPorUnidad
